# Applies the "0 Workshop.pptx" edits:
#  - Slide 1 (Text Box 2): "March 2021" -> "October 2021"
#  - Slide 3 (Rectangle 2): color the "R-stan practice 5: stanarm and model
#    checking in linear models." bullet a muted gray (White, Background 1,
#    Darker 35% == bg1 + lumMod 65%)
#  - Slide 3 (Rectangle 2): "Repeated measures (ppt)" ->
#    "Gaussian Linear models for repeated measures (ppt)"
#  - Add a (blank) reviewer comment by Juan Pedro steibel on slide 3

$p = $ppt.ActivePresentation

# --- Slide 1: update the workshop date -------------------------------------
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(2)
$dateRange = $dateShape.TextFrame.TextRange
$dateParagraph = $dateRange.Paragraphs(2, 1)
# Round-trip through an unrelated placeholder so the engine doesn't keep the
# shared "2021" suffix as a separate run - PowerPoint would just keep this
# as a single run with the original run formatting.
$dateParagraph.Text = "__TMP_DATE__"
$dateParagraph = $dateRange.Paragraphs(2, 1)
$dateParagraph.Text = "October 2021"

# --- Slide 3: bullet list updates -------------------------------------------
$slide3 = $p.Slides.Item(3)
$listShape = $slide3.Shapes.Item(2)
$listRange = $listShape.TextFrame.TextRange

# Paragraph 8 = "R-stan practice 5: stanarm and model checking in linear
# models." -> make it "White, Background 1, Darker 35%" (bg1 / lumMod 65%).
$practice5 = $listRange.Paragraphs(8, 1)
$practice5.Font.Color.RGB = 10921638

# Paragraph 11 = "Repeated measures (ppt)" -> "Gaussian Linear models for
# repeated measures (ppt)"
$repeated = $listRange.Paragraphs(11, 1)
$repeated.Text = "__TMP_REPEATED__"
$repeated = $listRange.Paragraphs(11, 1)
$repeated.Text = "Gaussian Linear models for repeated measures (ppt)"

# --- Slide 3: reviewer comment ----------------------------------------------
$comment = $slide3.Comments.Add(0.2943307086614173, 0.2638582677165354, "Juan Pedro steibel", "JPs", "")

Write-Host "edit complete"
